$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 1 (header row): add bold "Status + Passed/Failed? Bugs written up?" to column 3 ---
$r = $t.Cell(1, 3).Range
$r.InsertAfter("Status + ")
$r.Font.Bold = 1

$r = $t.Cell(1, 3).Range
$r.Collapse(0)
$r.InsertAfter("Passed/Failed? Bugs written up")
$r.Font.Bold = 1

$r = $t.Cell(1, 3).Range
$r.Collapse(0)
$r.InsertAfter("?")
$r.Font.Bold = 1

# --- "Position indicator (bottom right corner)" row: add "Connor" to column 2 ---
$r = $t.Cell(10, 2).Range
$r.InsertAfter("Connor")

# --- "Changing colors ..." row: add status note to column 3 ---
$r = $t.Cell(13, 3).Range
$r.InsertAfter("All passed")

$r = $t.Cell(13, 3).Range
$r.Collapse(0)
$r.InsertAfter(" – working on unit tests")
